$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 1-13) previously had no TR/EMA seed values - they are now
# explicit zeros (the rolling-ta "TR" warmup period written out as 0).
for ($i = 1; $i -le 13; $i++) {
    $ws.Cells.Item($i, 3).Value = 0
}

# C14 keeps its original seed formula (AVERAGE(B1:B14)) - untouched.

# From C15 down to C200 the EMA recursive formula
#   ((B{n} - C{n-1}) * (2/15)) + C{n-1}
# is replaced with a plain 14-period simple moving average:
#   AVERAGE(B{n-13}:B{n})
for ($i = 15; $i -le 200; $i++) {
    $ws.Cells.Item($i, 3).Formula = "=AVERAGE(B" + ($i - 13) + ":B" + $i + ")"
}

# Column widths: column C narrows and columns A & B now get explicit custom
# widths too (closest values this engine's 1/6-character rounding allows).
$ws.Columns.Item(1).ColumnWidth = 12.42578125 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 14.28515625 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 13.5703125 - (5/6)
